$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.43"
$ws.Range("E2").Value = "'-3.22%"
$ws.Range("D2:E2").Style = "Normal"
$ws.Range("D3").Value = "'31.39"
$ws.Range("E3").Value = "'-1.53%"
$ws.Range("D3:E3").Style = "Normal"
$ws.Range("D4").Value = "'4.961"
$ws.Range("E4").Value = "'-1.38%"
$ws.Range("D4:E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07368"
$ws.Range("E5").Value = "'-5.68%"
$ws.Range("D5:E5").Style = "Normal"
$ws.Range("D6").Value = "'1.819"
$ws.Range("E6").Value = "'-10.41%"
$ws.Range("D6:E6").Style = "Normal"
$ws.Range("D7").Value = "'7.662"
$ws.Range("E7").Value = "'-1.66%"
$ws.Range("D7:E7").Style = "Normal"
$ws.Range("D8").Value = "'3.748"
$ws.Range("E8").Value = "'-0.90%"
$ws.Range("D8:E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9101"
$ws.Range("E9").Value = "'-0.98%"
$ws.Range("D9:E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1645"
$ws.Range("E10").Value = "'-5.55%"
$ws.Range("D10:E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07615"
$ws.Range("E11").Value = "'-3.17%"
$ws.Range("D11:E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08182"
$ws.Range("E12").Value = "'-6.44%"
$ws.Range("D12:E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02989"
$ws.Range("E13").Value = "'-3.99%"
$ws.Range("D13:E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09955"
$ws.Range("E14").Value = "'-0.41%"
$ws.Range("D14:E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001501"
$ws.Range("E15").Value = "'-0.61%"
$ws.Range("D15:E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005641"
$ws.Range("E16").Value = "'-3.61%"
$ws.Range("D16:E16").Style = "Normal"
$ws.Range("D18").Value = "'3.467"
$ws.Range("E18").Value = "'0.19%"
$ws.Range("D18:E18").Style = "Normal"
$ws.Range("D19").Value = "'2.124"
$ws.Range("E19").Value = "'-6.24%"
$ws.Range("D19:E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3287"
$ws.Range("E20").Value = "'-0.18%"
$ws.Range("D20:E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1293"
$ws.Range("E21").Value = "'-0.02%"
$ws.Range("D21:E21").Style = "Normal"
$ws.Range("D22").Value = "'4.328"
$ws.Range("E22").Value = "'3.41%"
$ws.Range("D22:E22").Style = "Normal"
$ws.Range("E23").Value = "'9.17%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04496"
$ws.Range("E24").Value = "'-2.21%"
$ws.Range("D24:E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001226"
$ws.Range("E25").Value = "'-1.35%"
$ws.Range("D25:E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004051"
$ws.Range("E26").Value = "'-9.25%"
$ws.Range("D26:E26").Style = "Normal"
$ws.Range("E27").Value = "'0.08%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01638"
$ws.Range("E39").Value = "'-5.93%"
$ws.Range("D39:E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04383"
$ws.Range("E40").Value = "'-7.54%"
$ws.Range("D40:E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007452"
$ws.Range("E41").Value = "'4.11%"
$ws.Range("D41:E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1325"
$ws.Range("E42").Value = "'-2.20%"
$ws.Range("D42:E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002057"
$ws.Range("E43").Value = "'-1.16%"
$ws.Range("D43:E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01105"
$ws.Range("E44").Value = "'2.66%"
$ws.Range("D44:E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005990"
$ws.Range("E45").Value = "'-1.01%"
$ws.Range("D45:E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.04%"
$ws.Range("D46:E46").Style = "Normal"
$ws.Range("D47").Value = "'1.832"
$ws.Range("E47").Value = "'123.44%"
$ws.Range("D47:E47").Style = "Normal"
$ws.Range("E48").Value = "'-15.60%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("E50").Style = "Normal"
